$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @{
    "A1" = "Date"
    "B1" = "Na"
    "C1" = "K"
    "D1" = "Cl"
    "E1" = "ECO2"
    "F1" = "AGAP"
    "G1" = "AHOL"
    "H1" = "TBI"
    "I1" = "TP"
    "J1" = "GLOB"
    "K1" = "ALPI"
    "L1" = "TGL"
    "M1" = "CHOL"
    "N1" = "AST"
    "O1" = "ALTI"
    "P1" = "ALB"
    "Q1" = "A/G"
    "R1" = "GLUC"
    "S1" = "BUN"
    "T1" = "CA"
    "U1" = "CRE2"
    "V1" = "BN/CR"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}
